$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "Approuvée"
$ws.Range("G4").Value = "2025-05-23 09:04:01"
$ws.Range("H4").Value = "Magasinier"
$ws.Range("I4").Value = "Demande approuvée et stock mis à jour"

$ws.Range("A5").Value = "20250523_093023"
$ws.Range("B5").Value = "2025-05-23 09:30:23"
$ws.Range("C5").Value = "Elie"
$ws.Range("D5").Value = "{'chantier': 'Maintenance', 'urgence': 'Normal', 'date_souhaitee': '2025-05-23', 'produits': {'TS001': {'produit': 'Tournevis cruciforme', 'quantite': 2, 'emplacement': 'Atelier B'}, 'MH001': {'produit': 'Marteau 500g', 'quantite': 6, 'emplacement': 'Atelier B'}}}"
$ws.Range("E5").Value = "aa"
$ws.Range("F5").Value = "En attente"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Style = "Normal"
